{"js": "// The \"Scenarios\" bullet about a couple dining outdoors drops the\n// \"al fresco \" wording, and Word's auto-managed \"_GoBack\" bookmark\n// (which marks the location of the most recent edit) moves from the\n// end of the document to the spot where that edit just happened.\n\nconst document = context.document;\nconst body = document.body;\n\n// 1) Remove the old \"_GoBack\" bookmark wherever it currently sits\n//    (it starts out on the trailing empty Heading2 paragraph).\nconst goBack = document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) Trim \"al fresco \" out of the sentence, leaving the surrounding\n//    text (\"...enjoy dining \" + \"in their backyard...\") untouched.\nconst target = body.search(\"al fresco \", { matchCase: true, matchWholeWord: false });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error('Expected to find \"al fresco \" in the document body.');\n}\ntarget.items[0].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-create \"_GoBack\" right after \"enjoy dining \" \u2014 i.e. exactly\n//    where the text was just edited \u2014 matching Word's own behavior.\nconst editSpot = body.search(\"enjoy dining \", { matchCase: true, matchWholeWord: false });\neditSpot.load(\"items\");\nawait context.sync();\n\nif (editSpot.items.length === 0) {\n  throw new Error('Expected to find \"enjoy dining \" in the document body.');\n}\nconst insertionPoint = editSpot.items[0].getRange(Word.RangeLocation.after);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The \"Scenarios\" bullet about a couple dining outdoors drops the\n# \"al fresco \" wording, and Word's auto-managed \"_GoBack\" bookmark\n# (which marks the location of the most recent edit) moves from the\n# end of the document to the spot where that edit just happened.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark wherever it currently sits\n#    (it starts out on the trailing empty Heading2 paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Trim \"al fresco \" out of the sentence, leaving the surrounding\n#    text (\"...enjoy dining \" + \"in their backyard...\") untouched.\n$find = $d.Content\n[void]$find.Find.Execute(\"al fresco \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 3) Re-create \"_GoBack\" right after \"enjoy dining \" - i.e. exactly\n#    where the text was just edited - matching Word's own behavior.\n$mark = $d.Content\n[void]$mark.Find.Execute(\"enjoy dining \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$mark.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $mark)\n"}
